$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Change Percent" column (J) entirely, shrinking the used range.
$ws.Range("J1").EntireColumn.Delete()

# Drop the VTV holding (row 6) entirely; rows below shift up.
$ws.Range("A6").EntireRow.Delete()

# Reorder the F:I headers: Profit/Loss, Percentage Change, Todays Change,
# Todays Change Percent (previously Todays Change, Todays Change Percent,
# Profit/Loss, Percentage Change).
$ws.Range("F1").Value = "Profit/Loss"
$ws.Range("G1").Value = "Percentage Change"
$ws.Range("H1").Value = "Todays Change"
$ws.Range("I1").Value = "Todays Change Percent"

# Refresh every holding's numbers (rows 2-7), including the new
# F:I column order, and replace the final row's data with TSM.
$data = @(
    @("NVAX", 31.73, 844.58, 26798.5234, 8.659999847412109, -19484.46072887268, -72.70721762555276, -0.1100006103515625, -1.25),
    @("SMCI", 90.94, 360, 32738.4, 41.59999847412109, -17762.4005493164, -54.25555478983826, -0.4000015258789062, -0.95),
    @("NTNX", 73.16, 175, 12803, 65.34999847412109, -1366.750267028808, -10.67523445308762, -0.8000030517578125, -1.21),
    @("AVXL", 6.21, 2065, 12823.65, 8.970000267028809, 5699.40055141449, 44.44444874442526, -0.3499994277954102, -3.76),
    @("XLK", 234.6, 156.21, 36646.866, 235.9001007080078, 203.0887315979013, 0.5541776248967683, -0.049896240234375, -0.02),
    @("TSM", 188.22, 24.09, 4534.2198, 197.7149963378906, 228.7344617797852, 5.04462668042218, 3.31500244140625, 1.71)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $row++
}
